# Apply crypto price/volume updates (and two row re-orderings) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.461.46'
$ws.Range('D3').Value = '1.937.19'
$ws.Range('E3').Value = '  +4.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.39'
$ws.Range('E5').Value = '  +3.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4742'
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2868'
$ws.Range('E8').Value = '  +4.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06603'
$ws.Range('E9').Value = '  +4.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.51'
$ws.Range('E10').Value = '  +10.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '107.49'
$ws.Range('E11').Value = '  +27.52%  '
$ws.Range('D12').Value = '1.917.53'
$ws.Range('E12').Value = '  +3.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07594'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.133'
$ws.Range('E14').Value = '  +2.56%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6540'
$ws.Range('E15').Value = '  +4.95%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '310.88'
$ws.Range('E16').Value = '  +27.79%  '
$ws.Range('D17').Value = '30.483.75'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.03'
$ws.Range('E18').Value = '  +3.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007515'
$ws.Range('E20').Value = '  +2.64%  '
$ws.Range('D21').Value = '2.169.01'
$ws.Range('E21').Value = '  +3.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.200'
$ws.Range('E23').Value = '  +5.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.380'
$ws.Range('E24').Value = '  +7.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.312'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.46'
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.25'
$ws.Range('E27').Value = '  +13.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.061'
$ws.Range('E28').Value = '  +10.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1122'
$ws.Range('E29').Value = '  +10.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.348'
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.119'
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.956'
$ws.Range('E32').Value = '  +3.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05056'
$ws.Range('E33').Value = '  +4.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7438'
$ws.Range('E34').Value = '  +6.45%  '
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.717'
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01975'
$ws.Range('E37').Value = '  +4.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.711'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.035'
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8770'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.53'
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.903'
$ws.Range('E42').Value = '  +7.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9996'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '68.70'
$ws.Range('E44').Value = '  +10.36%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4170'
$ws.Range('E45').Value = '  +3.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.311'
$ws.Range('E46').Value = '  +2.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.337'
$ws.Range('E47').Value = '  +9.29%  '
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.79'
$ws.Range('E49').Value = '  +4.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05636'
$ws.Range('E50').Value = '  +1.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3853'
$ws.Range('E51').Value = '  +5.47%  '
